$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the header value from A1 into B1
$ws.Range("B1").Value = $ws.Range("A1").Value2

# Move the "Lime" translation from D2 to B2, and clear the old D2 cell
$ws.Range("B2").Value = $ws.Range("D2").Value2
$ws.Range("D2").ClearContents()
